$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''41.488.12'
$ws.Range("E2").Value = '  -2.64%  '

$ws.Range("D3").Value = '''2.468.00'
$ws.Range("E3").Value = '  -2.40%  '

$ws.Range("E4").Value = '  +0.79%  '

$ws.Range("D5").Value = '''311.95'
$ws.Range("E5").Value = '  -1.00%  '

$ws.Range("D6").Value = '''91.52'
$ws.Range("E6").Value = '  -6.90%  '

$ws.Range("D7").Value = '''0.540'
$ws.Range("E7").Value = '  -3.86%  '

$ws.Range("E8").Value = '  +0.66%  '

$ws.Range("E9").Value = '  -5.79%  '

$ws.Range("D10").Value = '''32.72'
$ws.Range("E10").Value = '  -7.15%  '

$ws.Range("E11").Value = '  -3.16%  '

$ws.Range("E12").Value = '  -0.35%  '

$ws.Range("D13").Value = '''2.851.66'
$ws.Range("E13").Value = '  -2.20%  '

$ws.Range("D14").Value = '''6.81'
$ws.Range("E14").Value = '  -5.65%  '

$ws.Range("D15").Value = '''2.469.87'
$ws.Range("E15").Value = '  -2.14%  '

$ws.Range("D16").Value = '''15.19'
$ws.Range("E16").Value = '  -0.05%  '

$ws.Range("D17").Value = '''0.776'
$ws.Range("E17").Value = '  -4.40%  '

$ws.Range("D18").Value = '''41.269.75'
$ws.Range("E18").Value = '  -3.19%  '

$ws.Range("D19").Value = '''6.24'
$ws.Range("E19").Value = '  -5.32%  '

$ws.Range("D20").Value = '''0.0₃0914'
$ws.Range("E20").Value = '  -2.76%  '

$ws.Range("D21").Value = '''70.50'
$ws.Range("E21").Value = '  +1.85%  '

$ws.Range("D22").Value = '''10.91'
$ws.Range("E22").Value = '  -9.92%  '

$ws.Range("D23").Value = '''234.83'
$ws.Range("E23").Value = '  -2.98%  '

$ws.Range("D24").Value = '''2.70'
$ws.Range("E24").Value = '  -5.35%  '

$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("E26").Value = '  -6.07%  '

$ws.Range("E27").Value = '  -6.18%  '

$ws.Range("E28").Value = '  -0.58%  '

$ws.Range("D29").Value = '''9.64'
$ws.Range("E29").Value = '  -3.60%  '

$ws.Range("D30").Value = '''35.63'
$ws.Range("E30").Value = '  -5.45%  '

$ws.Range("D31").Value = '''152.05'
$ws.Range("E31").Value = '  -2.32%  '

$ws.Range("E32").Value = '  -8.55%  '

$ws.Range("E33").Value = '  -4.27%  '

$ws.Range("D34").Value = '''2.56'
$ws.Range("E34").Value = '  -3.24%  '

$ws.Range("D35").Value = '''0.0750'
$ws.Range("E35").Value = '  -4.28%  '

$ws.Range("D36").Value = '''17.28'
$ws.Range("E36").Value = '  -1.17%  '

$ws.Range("D37").Value = '''2.96'
$ws.Range("E37").Value = '  -5.43%  '

$ws.Range("E38").Value = '  -7.30%  '

$ws.Range("D39").Value = '''0.112'
$ws.Range("E39").Value = '  -3.59%  '

$ws.Range("D40").Value = '''0.0993'
$ws.Range("E40").Value = '  -8.09%  '

$ws.Range("D41").Value = '''4.00'
$ws.Range("E41").Value = '  -5.83%  '

$ws.Range("E42").Value = '  +0.99%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '''19.25'
$ws.Range("E43").Value = '  -10.08%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '''1.959.10'
$ws.Range("E44").Value = '  -3.39%  '

$ws.Range("D45").Value = '''0.0280'
$ws.Range("E45").Value = '  -5.19%  '

$ws.Range("D46").Value = '''2.92'
$ws.Range("E46").Value = '  -8.95%  '

$ws.Range("D47").Value = '''8.58'
$ws.Range("E47").Value = '  -3.27%  '

$ws.Range("D48").Value = '''2.719.88'
$ws.Range("E48").Value = '  -1.75%  '

$ws.Range("D49").Value = '''95.26'
$ws.Range("E49").Value = '  -4.80%  '

$ws.Range("D50").Value = '''67.52'
$ws.Range("E50").Value = '  -6.04%  '

$ws.Range("E51").Value = '  -7.33%  '
